$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H19").Value = 811.7143
$ws.Range("I19").Value = 315.2
$ws.Range("J19").Value = 1263.091
$ws.Range("K19").Value = 315.2
$ws.Range("L19").Value = 1263.091
$ws.Range("M19").Value = -140.2
$ws.Range("N19").Value = -1613.091
$ws.Range("H29").Value = 1877.2
$ws.Range("I29").Value = 734.3333
$ws.Range("J29").Value = 2367
$ws.Range("K29").Value = 2202.9999
$ws.Range("L29").Value = 7101
$ws.Range("M29").Value = -1921.9999
$ws.Range("N29").Value = -7663
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H112").Value = 2513.4062
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 2559
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 7677
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -9893
$ws.Range("H138").Value = 2197.6326
$ws.Range("I138").Value = 1092.4286
$ws.Range("J138").Value = 2282.6484
$ws.Range("K138").Value = 3277.2858
$ws.Range("L138").Value = 6847.9452
$ws.Range("M138").Value = 1862.7142
$ws.Range("N138").Value = -17127.9452
$ws.Range("H141").Value = 13086.111
$ws.Range("I141").Value = 14221.875
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 42665.625
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -37485.625
$ws.Range("N141").Value = -22360

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 2921.1968
$ws.Range("I32").Value = 3014.1924
$ws.Range("K32").Value = 3014.1924
$ws.Range("M32").Value = -2727.1924
$ws.Range("H61").Value = 952.88
$ws.Range("I61").Value = 659.75
$ws.Range("K61").Value = 659.75
$ws.Range("M61").Value = -447.75
$ws.Range("H136").Value = 952.88
$ws.Range("I136").Value = 659.75
$ws.Range("K136").Value = 1979.25
$ws.Range("M136").Value = 570.75

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H22").Value = 151.33333
$ws.Range("I22").Value = 141.6
$ws.Range("K22").Value = 141.6
$ws.Range("M22").Value = 31.40000000000001

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 1332.1428
$ws.Range("I31").Value = 1072.5385
$ws.Range("J31").Value = 1754
$ws.Range("K31").Value = 1072.5385
$ws.Range("L31").Value = 1754
$ws.Range("M31").Value = -777.5385000000001
$ws.Range("N31").Value = -2344
$ws.Range("H34").Value = 1332.1428
$ws.Range("I34").Value = 1072.5385
$ws.Range("J34").Value = 1754
$ws.Range("K34").Value = 1072.5385
$ws.Range("L34").Value = 1754
$ws.Range("M34").Value = -870.5385000000001
$ws.Range("N34").Value = -2158
$ws.Range("H58").Value = 1308.5555
$ws.Range("I58").Value = 1296.1666
$ws.Range("K58").Value = 1296.1666
$ws.Range("M58").Value = -1093.1666
$ws.Range("H134").Value = 10418097
$ws.Range("I134").Value = 13334716
$ws.Range("K134").Value = 40004148
$ws.Range("M134").Value = -40001613
$ws.Range("H136").Value = 1308.5555
$ws.Range("I136").Value = 1296.1666
$ws.Range("K136").Value = 3888.4998
$ws.Range("M136").Value = -1338.4998

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H80").Value = 2159.8
$ws.Range("I80").Value = 1199.75
$ws.Range("K80").Value = 3599.25
$ws.Range("M80").Value = -2663.25
$ws.Range("H81").Value = 2628.6086
$ws.Range("J81").Value = 2852.5
$ws.Range("L81").Value = 8557.5
$ws.Range("N81").Value = -10803.5
$ws.Range("H83").Value = 2159.8
$ws.Range("I83").Value = 1199.75
$ws.Range("K83").Value = 10797.75
$ws.Range("M83").Value = -6117.75
$ws.Range("H84").Value = 2628.6086
$ws.Range("J84").Value = 2852.5
$ws.Range("L84").Value = 25672.5
$ws.Range("N84").Value = -36904.5
$ws.Range("H88").Value = 3499.0625
$ws.Range("J88").Value = 7664.1665
$ws.Range("L88").Value = 22992.4995
$ws.Range("N88").Value = -23848.4995
$ws.Range("H91").Value = 3499.0625
$ws.Range("J91").Value = 7664.1665
$ws.Range("L91").Value = 22992.4995
$ws.Range("N91").Value = -25956.4995
$ws.Range("H107").Value = 3878.4062
$ws.Range("I107").Value = 647.9524
$ws.Range("K107").Value = 1943.8572
$ws.Range("M107").Value = -23.85719999999992
$ws.Range("H113").Value = 519.0213
$ws.Range("I113").Value = 440.76923
$ws.Range("J113").Value = 548.94116
$ws.Range("K113").Value = 1322.30769
$ws.Range("L113").Value = 1646.82348
$ws.Range("M113").Value = 847.6923099999999
$ws.Range("N113").Value = -5986.82348
$ws.Range("H122").Value = 842.75
$ws.Range("I122").Value = 548.6
$ws.Range("J122").Value = 1333
$ws.Range("K122").Value = 4937.400000000001
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -2487.400000000001
$ws.Range("N122").Value = -16897
$ws.Range("H124").Value = 2312.5
$ws.Range("I124").Value = 666.6667
$ws.Range("K124").Value = 2000.0001
$ws.Range("M124").Value = 2909.9999
$ws.Range("H131").Value = 20834760
$ws.Range("I131").Value = 142857800
$ws.Range("J131").Value = 1560.9512
$ws.Range("K131").Value = 428573400
$ws.Range("L131").Value = 4682.8536
$ws.Range("M131").Value = -428568360
$ws.Range("N131").Value = -14762.8536
$ws.Range("H137").Value = 11113.5
$ws.Range("I137").Value = 2407.5
$ws.Range("J137").Value = 15466.5
$ws.Range("K137").Value = 7222.5
$ws.Range("L137").Value = 46399.5
$ws.Range("M137").Value = -2122.5
$ws.Range("N137").Value = -56599.5

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H39").Value = 20100
$ws.Range("J39").Value = 20100
$ws.Range("L39").Value = 20100
$ws.Range("N39").Value = -21164
$ws.Range("H126").Value = 2193.457
$ws.Range("I126").Value = 1354.2
$ws.Range("J126").Value = 4291.6
$ws.Range("K126").Value = 4062.6
$ws.Range("L126").Value = 12874.8
$ws.Range("M126").Value = -1592.6
$ws.Range("N126").Value = -17814.8

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
